$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AddJobTitles")
$ws2 = $wb.Worksheets.Item("DeleteJobTitles")

# On "AddJobTitles": the 4th job title is edited from "TESTER" to "ALFA",
# then the cursor moves down to A5 (as when typing Enter after an edit).
$ws1.Activate()
$ws1.Range("A4").Value = "ALFA"
[void]$ws1.Range("A5").Select()

# On "DeleteJobTitles": "ALFA" is added to A3 (previously empty),
# the sheet becomes the active tab, and the selection ends on E4.
$ws2.Activate()
$ws2.Range("A3").Value = "ALFA"
$ws2.PageSetup.Orientation = 1
[void]$ws2.Range("E4").Select()
